$wb = $excel.ActiveWorkbook

# --- Rename the header cells on the two existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$sheetCount = $wb.Worksheets.Count
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Copy the existing bold/centered header style (A1 on "Weekly Quantity") onto the
# new header row so the new sheet's formatting matches the workbook's convention.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$rows = @(
    @(44934.99999999999, 867, -707.6262833788974, 2427.541732356647),
    @(44941.99999999999, 869, -502.3718782367143, 2323.550635821488),
    @(44955.99999999999, 873, -687.7315539440028, 2385.172307721542),
    @(44962.99999999999, 875, -688.0635249390308, 2336.481957694816),
    @(44969.99999999999, 877, -582.059403644872, 2311.354301982832),
    @(44976.99999999999, 879, -716.1997314805416, 2350.521377338677),
    @(44983.99999999999, 881, -549.5780533309635, 2373.451419620368),
    @(45018.99999999999, 892, -666.1446546216977, 2222.083928813098),
    @(45032.99999999999, 896, -785.1845166809987, 2323.412111660037),
    @(45039.99999999999, 898, -628.1359634153733, 2409.663423724364),
    @(45046.99999999999, 900, -621.7010980443848, 2389.332033001493),
    @(45053.99999999999, 902, -611.8099400913554, 2475.645416138929),
    @(45060.99999999999, 904, -599.2299612059238, 2433.301752433723),
    @(45067.99999999999, 907, -475.8760193182584, 2514.291465587828),
    @(45074.99999999999, 909, -733.9157492568071, 2313.035574259691),
    @(45095.99999999999, 915, -613.6928300760612, 2435.025858482564),
    @(45102.99999999999, 917, -519.481379924644, 2391.126552915374),
    @(45137.99999999999, 927, -659.4149138612594, 2434.140241075021),
    @(45144.99999999999, 930, -516.7953405593296, 2452.728960161011),
    @(45151.99999999999, 932, -569.1777534129652, 2595.792719138955),
    @(45158.99999999999, 934, -616.6433196867523, 2446.730933170715),
    @(45186.99999999999, 942, -630.987039287863, 2515.60480577079),
    @(45193.99999999999, 944, -543.3840583671979, 2317.809206215936),
    @(45200.99999999999, 946, -472.3726241981212, 2449.657963142407),
    @(45221.99999999999, 953, -583.487773017754, 2441.325994966036),
    @(45228.99999999999, 955, -607.1253863292711, 2318.777787259204),
    @(45235.99999999999, 957, -591.436135923776, 2461.130193378784),
    @(45242.99999999999, 959, -420.1282069987744, 2556.652082147255),
    @(45249.99999999999, 961, -482.187146429156, 2362.505736787839),
    @(45256.99999999999, 963, -470.2396855263331, 2396.679087644793),
    @(45270.99999999999, 967, -542.2480313666383, 2492.218874672741),
    @(45277.99999999999, 969, -560.7950380487877, 2390.934504793229),
    @(45298.99999999999, 976, -480.5670683028305, 2551.452865937633),
    @(45305.99999999999, 978, -581.9542379565775, 2529.679595170501),
    @(45326.99999999999, 984, -543.0949971655706, 2418.659597774833),
    @(45361.99999999999, 994, -472.5286679456196, 2484.216536915694),
    @(45368.99999999999, 997, -476.6430973054055, 2600.348401644927),
    @(45375.99999999999, 999, -413.3443639552619, 2658.315331624307),
    @(45382.99999999999, 1001, -517.8179341760381, 2486.070438161166),
    @(45389.99999999999, 1003, -563.4928582754717, 2384.844417611644),
    @(45403.99999999999, 1007, -528.0284369691774, 2518.637269974447),
    @(45410.99999999999, 1009, -513.3205641072523, 2531.933811446925),
    @(45417.99999999999, 1011, -520.8308896475933, 2517.388903939485),
    @(45424.99999999999, 1013, -476.4606472390046, 2550.406455056086),
    @(45438.99999999999, 1018, -332.3003757074608, 2462.08558974142),
    @(45445.99999999999, 1020, -514.0090382393123, 2465.791891845433),
    @(45452.99999999999, 1022, -475.6608074526677, 2599.085961040168),
    @(45459.99999999999, 1024, -533.0046637000653, 2486.037953335565),
    @(45466.99999999999, 1026, -490.1148924918278, 2526.034798397764),
    @(45494.99999999999, 1034, -666.3493775725404, 2519.918462055786),
    @(45529.99999999999, 1045, -486.308550901864, 2467.0016232671),
    @(45536.99999999999, 1047, -468.9386375287517, 2541.21083860675),
    @(45543.99999999999, 1049, -591.8998733272168, 2539.616940592382),
    @(45550.99999999999, 1051, -509.8565236996491, 2534.946205482458),
    @(45557.99999999999, 1053, -615.3555478491078, 2430.999671194031),
    @(45564.99999999999, 1055, -499.4983932298556, 2465.224031543906),
    @(45571.99999999999, 1057, -520.2818609085949, 2580.54758040737),
    @(45578.99999999999, 1059, -539.7761164292399, 2524.162311553299),
    @(45585.99999999999, 1061, -426.932637005753, 2523.340233989008),
    @(45592.99999999999, 1064, -476.3431771330168, 2557.403970465115),
    @(45599.99999999999, 1066, -600.3376218264655, 2485.859400736872),
    @(45606.99999999999, 1068, -509.1127437837943, 2513.678809639368),
    @(45620.99999999999, 1072, -372.9444926712934, 2564.928199125631),
    @(45627.99999999999, 1074, -446.3245673426138, 2663.313145972819),
    @(45634.99999999999, 1076, -399.241892941304, 2627.301189820449),
    @(45641.99999999999, 1078, -431.9015077057291, 2609.660871567822),
    @(45648.99999999999, 1080, -413.2830631657171, 2687.517182178803),
    @(45655.99999999999, 1082, -343.4950577141801, 2540.64295832068),
    @(45662.99999999999, 1085, -337.1848152018247, 2521.885631320898),
    @(45669.99999999999, 1087, -348.9419208029659, 2514.52446668133),
    @(45676.99999999999, 1089, -506.8560575507939, 2620.583114984674)
)

$rowIndex = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($rowIndex, 1).Value = $row[0]
    $wsForecast.Cells.Item($rowIndex, 2).Value = $row[1]
    $wsForecast.Cells.Item($rowIndex, 3).Value = $row[2]
    $wsForecast.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex = $rowIndex + 1
}

# Copy the existing date-formatted style (A2 on "Weekly Quantity") onto the whole
# "ds" column of data rows so it renders/serializes the same way as the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A72").PasteSpecial(-4122)
